$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BookedHours")

# --- Week 2 booked hours entry (row 4) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "29/07/2025"
$ws.Range("C4").Value = 0.50694444444444442
$ws.Range("D4").Value = "29/07/2025"
$ws.Range("D4").NumberFormat = $ws.Range("C4").NumberFormat
$ws.Range("E4").Value = 0.53819444444444442

# --- Row 5 stray ":" entry in the Begin-Time column ---
$ws.Range("C5").Value = ":"

$ws.Range("F4").Value = "Group meeting"
$ws.Range("G4").Value = 0.75
$ws.Range("H4").Value = 0.75

# --- Update the title / instructions banner in A1 ---
$c = $ws.Range("A1")
$c.Value = "CITS3200 Project Billed Hours Record for Xuan Truong Nguyen and end date and hour, plus a brief description of the activity. At the end of each week send a copy of the sheet as it currently stands to your group's manager for recording on the group TimeSheet. Date in the form DD/MM/YY and time as XX:YY (24hr clock)"
$len = $c.Value.Length
$startSecond = "CITS3200 Project Billed Hours Record for Xuan Truong Nguyen".Length + 1
$chars = $c.Characters($startSecond, $len - $startSecond + 1)
$chars.Font.Bold = $false
$chars.Font.Size = 12
$chars.Font.Name = "Lucida Sans"
$chars.Font.ColorIndex = -4105

# --- View state: selected cell and zoom level ---
$ws.Activate() | Out-Null
$ws.Range("F9").Select() | Out-Null
$excel.ActiveWindow.Zoom = 82
